$d = $word.ActiveDocument

function Set-ParaText($index, $level, $text) {
    $p = $d.Paragraphs.Item($index)
    $r = $d.Range($p.Range.Start, $p.Range.End)
    $r.Text = $text
    $d.Paragraphs.Item($index).Range.ListFormat.ListLevelNumber = $level
}

function Insert-ParaAfter($index, $level, $text) {
    $p = $d.Paragraphs.Item($index)
    $r = $d.Range($p.Range.Start, $p.Range.End)
    $r.InsertParagraphAfter()
    $newP = $d.Paragraphs.Item($index + 1)
    $newP.Range.Text = $text
    $d.Paragraphs.Item($index + 1).Range.ListFormat.ListLevelNumber = $level
}

# --- Simple text-only fixes (run merges elsewhere in the diff do not change
#     visible text, so they are skipped) -------------------------------------

# Heading: "Web Interface Requirements" -> "Web Interface/Routing Processor Requirements"
$d.Content.Find.Execute("Web Interface Requirements", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Web Interface/Routing Processor Requirements", 2) | Out-Null

# "edges connecting nodes" -> "neighbor edges" (appears 4 times across the
# "hide/remove/add" GUI-element requirements plus the default-display one)
$d.Content.Find.Execute("edges connecting nodes", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "neighbor edges", 2) | Out-Null

# --- Restructure the "When the Web UI is initialized ... shall:" sub-list ---
# Reuse the 8 existing sub-list paragraphs (old indices 14-21) in place,
# rewriting their text/list-level to match the new content order.
Set-ParaText 14 2 "Calculate single-path routing tables for every source, destination pair in the network."
Set-ParaText 15 2 "Update the routing table for every node in the network."
Set-ParaText 16 2 "Send the updated network state (nodes, neighbor edges, routes) to the Web UI."
Set-ParaText 17 1 "The routing processor shall check every network element removal request to prevent creation of orphan nodes."
Set-ParaText 18 1 "When a user deletes a neighbor edge in the Web UI, the routing processor shall:"
Set-ParaText 19 2 "Check the routing impact of deleting this edge by calculating the number of affected of routes."
Set-ParaText 20 3 "If there are no routes using this neighbor edge, then the routing processor shall proceed with edge deletion."
Set-ParaText 21 3 "If there are 1 or more routes using this neighbor edge, the routing processor shall attempt to recalculate replacement routes which do not include the deleted neighbor edge."

# Insert three brand-new paragraphs after old paragraph 21.
Insert-ParaAfter 21 4 "If there are no routes that can replace the affected routes, the routing processor shall halt the edge deletion process, request that the Web UI prompt the user for confirmation, and then continue with edge deletion once the operation is confirmed by the user."
Insert-ParaAfter 22 4 "If there are routes that can replace the affected routes, the routing processor shall update the relevant node routing tables and then proceed with edge deletion."
Insert-ParaAfter 23 2 "Remove the nodes connected by the edge from node neighbor tables."

# The original last paragraph (bookmark _GoBack lives here) shifted from
# index 22 to index 25; give it its new wording while keeping the bookmark.
Set-ParaText 25 2 "Remove any routes present in all routing tables that uses this edge for multi-hop communication."
